$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateVal = Get-Date -Year 2025 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0

# New column header "Trening" - copy formatting (bold + border) from existing header cell
$ws.Range("F1").Value = "Trening"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats

# Row 2: convert A2 from text date "27.01.2025" to a real datetime value
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").Value = $dateVal
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F2").Value = "Duża Gra"

# Row 3: same conversion
$ws.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3").Value = $dateVal
$ws.Range("F3").Value = "Duża Gra"

# New row 4 - split training part ("Mała Gra")
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A4").Value = $dateVal
$ws.Range("B4:D4").Font.Bold = $false
$ws.Range("E4").Value = "10-15"
$ws.Range("F4").Value = "Mała Gra"

# New row 5
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A5").Value = $dateVal
$ws.Range("B5:D5").Font.Bold = $false
$ws.Range("E5").Value = "5-10"
$ws.Range("F5").Value = "Mała Gra"
